{"js": "// Update the \"Strengths\" section of the resume:\n//  1. Re-order the language list sentence so C# leads (with its new\n//     parenthetical) and Python/Ruby move down the list.\n//  2. Add a brand-new paragraph right after it describing familiarity\n//     with parallel/concurrent computation models.\n\nconst oldSentence =\n  \"Including, but not limit to: Python, Ruby, C++ (including C++11 features), \" +\n  \"C#, Java, Actionscipt 3 (and its assembly), Objective-C (without Cocoa), \" +\n  \"Erlang, Common Lisp.\";\n\nconst newSentence =\n  \"Including, but not limit to: C# (including TPL, async/await, dynamic IL \" +\n  \"generation), C++ (including C++11 features), Python, Ruby, Java, \" +\n  \"Actionscipt 3 (and its assembly), Objective-C (without Cocoa), Erlang, \" +\n  \"Common Lisp.\";\n\nconst newParagraphText =\n  \"Proficient or familiar with parallel/concurrent computation models: \" +\n  \"C#/TPL, Erlang/Actor, Python/gevent.\";\n\nconst body = context.document.body;\n\n// Locate the exact sentence as a Range so we touch only that run's text,\n// leaving the preceding \"Proficient or familiar...\" run / <w:cr/> intact.\nconst matches = body.search(oldSentence, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not find the languages sentence to update.\");\n}\n\nconst target = matches.items[0];\n\n// Swap in the reordered sentence.\ntarget.insertText(newSentence, \"Replace\");\n\n// Insert the new \"parallel/concurrent computation models\" paragraph\n// immediately after the paragraph that holds the sentence we just edited.\nconst containingParagraph = target.paragraphs.getFirst();\ncontainingParagraph.insertParagraph(newParagraphText, \"After\");\n\nawait context.sync();\n", "ps1": "# Update the \"Strengths\" section of the resume:\n#  1. Re-order the language list sentence so C# leads (with its new\n#     parenthetical) and Python/Ruby move down the list.\n#  2. Add a brand-new paragraph right after it describing familiarity\n#     with parallel/concurrent computation models.\n\n$d = $word.ActiveDocument\n\n$oldSentence = \"Including, but not limit to: Python, Ruby, C++ (including C++11 features), C#, Java, Actionscipt 3 (and its assembly), Objective-C (without Cocoa), Erlang, Common Lisp.\"\n$newSentence = \"Including, but not limit to: C# (including TPL, async/await, dynamic IL generation), C++ (including C++11 features), Python, Ruby, Java, Actionscipt 3 (and its assembly), Objective-C (without Cocoa), Erlang, Common Lisp.\"\n$newParagraphText = \"Proficient or familiar with parallel/concurrent computation models: C#/TPL, Erlang/Actor, Python/gevent.\"\n\n# Find the paragraph that holds the languages sentence.\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Including, but not limit to:*\") {\n        $targetParagraph = $p\n        break\n    }\n}\nif ($null -eq $targetParagraph) {\n    throw \"Could not find the languages paragraph to update.\"\n}\n\n# Replace just the sentence text (leaves the preceding run / manual line\n# break untouched) within that paragraph's range.\n$findRange = $targetParagraph.Range.Duplicate\n$findRange.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null\n\n# Insert a brand new paragraph right after it with the new sentence about\n# parallel/concurrent computation models.\n$targetParagraph.Range.InsertParagraphAfter() | Out-Null\n$newParagraph = $targetParagraph.Next()\n$newParagraph.Range.Text = $newParagraphText\n"}
